# Add "availability" and "staged_products" columns to the Resource sheet,
# inserted right before the existing "capacity" column (which shifts,
# along with "skills", "aggregates" and "kwargs", two columns to the right).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at I:J (pushing capacity/skills/aggregates/kwargs
# from I:L to K:N). Inserting whole columns preserves the existing header
# formatting (bold font + border, style index 1) for the shifted cells and
# for the newly created cells in the inserted range.
$ws.Range("I1:J1").EntireColumn.Insert()

# Populate the two newly inserted header cells.
$ws.Cells.Item(1, 9).Value = "availability"
$ws.Cells.Item(1, 10).Value = "staged_products"
